$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 3 ("Headshot in header to the left of name and email"),
# shifting all subsequent rows up by one.
$ws.Rows.Item(3).Delete()

# Reflect the post-edit active selection on Sheet1.
$ws.Range("A5").Select()
